# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Updates the "Estado de Cuenta" worksheet:
#  - Refreshes header totals (Valor Mora, Cant. Trabajadores, Cant. Periodos)
#  - Extends worker 8755981 (MARIO DE JESUS ZAMBRANO PEDRAZA) history with
#    one additional period (2508) and re-levels the "Valor Mora" amounts so
#    the most recent five periods (2208-2212) carry 40000 and the remaining
#    ones (2301-2508) carry 35112
#  - Adds three new workers (first new batch of accounts) each with a single
#    2508 period row
#  - Keeps the trailing signature block intact, shifted down to make room
#    for the new rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room for the 4 additional data rows (new row 52 for the existing
#    worker's extra period, plus 3 new rows for the new workers). This
#    shifts the trailing signature block (old rows 56-57) down to 60-61,
#    and the merged cells / sheet dimension follow automatically.
# ---------------------------------------------------------------------
$ws.Rows("52:55").Insert()

# ---------------------------------------------------------------------
# 2. Header summary values
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 1508664   # VALOR MORA (total)
$ws.Range("C13").Value = 4         # Cant. Trabajadores
$ws.Range("F13").Value = 37        # Cant. Periodos

# ---------------------------------------------------------------------
# 3. Re-level "Valor Mora" (column F) for worker 8755981's period rows
#    (rows 16-20 are now the newest periods 2208-2212 -> 40000, the rest,
#    rows 21-51, stay/become 35112).
# ---------------------------------------------------------------------
foreach ($r in 16..20) {
    $ws.Cells.Item($r, 6).Value = 40000
}
foreach ($r in 21..51) {
    $ws.Cells.Item($r, 6).Value = 35112
}

# ---------------------------------------------------------------------
# 4. Periods (column E) for worker 8755981 now run ascending 2208 -> 2508
#    across rows 16-52 (37 rows total).
# ---------------------------------------------------------------------
$periods = @( `
    "2208","2209","2210","2211","2212", `
    "2301","2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312", `
    "2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412", `
    "2501","2502","2503","2504","2505","2506","2507","2508" `
)
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "8755981"
    $ws.Cells.Item($r, 4).Value = "MARIO DE JESUS ZAMBRANO PEDRAZA"
    $ws.Cells.Item($r, 5).Value = $periods[$i]
    $ws.Cells.Item($r, 7).Value = 877803
}
# Row 52 is the new row carrying period 2508 for the same worker
$ws.Cells.Item(52, 6).Value = 35112

# ---------------------------------------------------------------------
# 5. New workers (first batch of new accounts), each a single 2508 row
# ---------------------------------------------------------------------
$newWorkers = @( `
    @{ Row = 53; Doc = "73129022";   Name = "FRANKLIN JOSE MARTINEZ GUTIERREZ"; Mora = 71200; Salario = 1780000 }, `
    @{ Row = 54; Doc = "1050945965"; Name = "MOISES FONTALVO GARCIA";           Mora = 56940; Salario = 1423500 }, `
    @{ Row = 55; Doc = "1002060363"; Name = "SAMUEL RAFAEL PEREIRA VASQUEZ";    Mora = 56940; Salario = 1423500 } `
)

foreach ($w in $newWorkers) {
    $r = $w.Row
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = $w.Doc
    $ws.Cells.Item($r, 4).Value = $w.Name
    $ws.Cells.Item($r, 5).Value = "2508"
    $ws.Cells.Item($r, 6).Value = $w.Mora
    $ws.Cells.Item($r, 7).Value = $w.Salario
}

# ---------------------------------------------------------------------
# 6. Formatting: give the new data rows (51-55) the same bordered look as
#    the rest of the worker table, and keep the closing (thicker bottom
#    border) style on the new last data row (55).
#    NOTE: source/destination ranges must be the same size, otherwise
#    PasteSpecial tiles the copied block repeatedly into the following
#    rows.
# ---------------------------------------------------------------------
$ws.Range("B16:J19").Copy()
$ws.Range("B51:J54").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# row 51 previously carried the table's closing (thicker bottom) border;
# row 55 is now the closing row, so give it that same look.
$ws.Range("B20:J20").Copy()
$ws.Range("B55:J55").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

foreach ($col in @("B","C","D","E","F","G","H","I","J")) {
    $cell = $ws.Range($col + "55")
    $cell.Borders.LineStyle = 1
    $cell.Borders.Weight = 2
}

Write-Host "Estado de cuenta actualizado: base de datos EC refrescada y parte 1 de nuevos estados de cuenta agregada."
